$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.178.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.764.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.763.57'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('E10').Value = '  +3.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.41'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E14').Value = '  +2.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.390.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.757.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.223.07'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +19.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '494.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.730'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  +8.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.65'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.906.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.699.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.28%  '
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.137'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.326'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('E42').Value = '  +5.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '430.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.67'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('E46').Value = '  +1.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.795.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0354'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.83%  '
